# 2017 ppt fixes. Added course outlines. Updated WhatsNew
#
# Slide 1 (title slide): the subtitle placeholder currently reads
# "FME 2016 Training" and needs to become "FME 2017" - as if the author
# selected "2016 Training" and typed "2017" over it, leaving the
# leading "FME " run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder robustly (ppPlaceholderSubtitle = 4)
$subtitle = $null
for ($i = 1; $i -le $s.Shapes.Placeholders.Count; $i++) {
    $ph = $s.Shapes.Placeholders.Item($i)
    if ($ph.PlaceholderFormat.Type -eq 4) {
        $subtitle = $ph
        break
    }
}
if ($subtitle -eq $null) {
    # fall back to the known shape if placeholder lookup fails
    $subtitle = $s.Shapes.Item(2)
}

$tr = $subtitle.TextFrame.TextRange
$full = $tr.Text
$prefix = "FME "

if ($full.StartsWith($prefix)) {
    # Select everything after "FME " (i.e. "2016 Training") and type
    # "2017" over it, the same way the author edited the slide.
    $startPos = $prefix.Length + 1
    $selLen = $tr.Length - $prefix.Length
    $target = $tr.Characters($startPos, $selLen)
    $target.Text = "2017"
} else {
    # Fallback for unexpected existing content.
    $tr.Text = "FME 2017"
}
